# Before using the code edit the config file
#
# The "status" / "Pass" block that already lives in columns E:J is being
# replicated across 11 more columns (K:U) so every credential row reports
# its status once per config entry. Column J (the right-most existing
# status column) is copied across for both content and formatting so the
# new header cells pick up the same colored fill the other status
# headers (E1:J1) use.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCols = @("K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

foreach ($col in $newCols) {
    # Copy the whole J column (header + data) into place.
    $ws.Range("J1:J6").Copy()
    $ws.Range($col + "1").PasteSpecial()

    # Re-apply J1's cell formatting (its colored fill) onto the new header
    # cell - PasteSpecial(xlPasteAll) above only carries the value/text.
    $ws.Range("J1").Copy()
    $ws.Range($col + "1").PasteSpecial(-4122)

    $ws.Columns($col).ColumnWidth = 5.518229166666667
}

$excel.CutCopyMode = 0
